# Update NATMI TPM-derived metrics (ligand/receptor expression & edge specificity
# columns) with newly recomputed TPM values, per "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1212753333333333
$ws.Range("H2").Value = 0.363826
$ws.Range("M2").Value = 0.665744
$ws.Range("N2").Value = 1.997232
$ws.Range("O2").Value = 0.07058985944777574
$ws.Range("P2").Value = 0.07058985944777574
$ws.Range("Q2").Value = 0.08073832551466667
$ws.Range("R2").Value = 0.7266449296319999
$ws.Range("S2").Value = 0.07058985944777574
$ws.Range("T2").Value = 0.07058985944777574

# Row 3
$ws.Range("G3").Value = 0.1212753333333333
$ws.Range("H3").Value = 0.363826
$ws.Range("O3").Value = 0.1615251915556908
$ws.Range("P3").Value = 0.1615251915556908
$ws.Range("Q3").Value = 0.1847471236897777
$ws.Range("R3").Value = 1.662724113208
$ws.Range("S3").Value = 0.1615251915556908
$ws.Range("T3").Value = 0.1615251915556908

# Row 4
$ws.Range("G4").Value = 0.1212753333333333
$ws.Range("H4").Value = 0.363826
$ws.Range("M4").Value = 1.786190333333334
$ws.Range("N4").Value = 5.358571
$ws.Range("O4").Value = 0.1893925060938975
$ws.Range("P4").Value = 0.1893925060938975
$ws.Range("Q4").Value = 0.2166208280717778
$ws.Range("R4").Value = 1.949587452646
$ws.Range("S4").Value = 0.1893925060938975
$ws.Range("T4").Value = 0.1893925060938975

# Row 5
$ws.Range("G5").Value = 0.1212753333333333
$ws.Range("H5").Value = 0.363826
$ws.Range("M5").Value = 2.606414666666667
$ws.Range("N5").Value = 7.819244
$ws.Range("O5").Value = 0.2763621526932594
$ws.Range("P5").Value = 0.2763621526932594
$ws.Range("Q5").Value = 0.3160938075048889
$ws.Range("R5").Value = 2.844844267544
$ws.Range("S5").Value = 0.2763621526932594
$ws.Range("T5").Value = 0.2763621526932594

# Row 6
$ws.Range("G6").Value = 0.1212753333333333
$ws.Range("H6").Value = 0.363826
$ws.Range("M6").Value = 1.713409666666666
$ws.Range("N6").Value = 5.140229
$ws.Range("O6").Value = 0.1816754601565471
$ws.Range("P6").Value = 0.1816754601565471
$ws.Range("Q6").Value = 0.2077943284615555
$ws.Range("R6").Value = 1.870148956154
$ws.Range("S6").Value = 0.1816754601565471
$ws.Range("T6").Value = 0.1816754601565471

# Row 7
$ws.Range("G7").Value = 0.1212753333333333
$ws.Range("H7").Value = 0.363826
$ws.Range("M7").Value = 1.136028333333333
$ws.Range("N7").Value = 3.408085
$ws.Range("O7").Value = 0.1204548300528295
$ws.Range("P7").Value = 0.1204548300528295
$ws.Range("Q7").Value = 0.1377722148011111
$ws.Range("R7").Value = 1.23994993321
$ws.Range("S7").Value = 0.1204548300528295
$ws.Range("T7").Value = 0.1204548300528295
